# "Generación del Estudio DT"
# - Adds 3 new "Trabajo de Grado" course rows to PlanEstudios (rows 66-68)
# - Updates the Libre Elección credit total on Agrupaciones (D11: 39 -> 28)
# - Clears the (redundant / unused) explicit cell styles that Excel had
#   been carrying on a bunch of cells across all three sheets
# - Leaves Agrupaciones as the active/selected sheet, with PlanEstudios
#   scrolled down near the newly appended rows

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PlanEstudios")
$ws2 = $wb.Worksheets.Item("Agrupaciones")
$ws3 = $wb.Worksheets.Item("Equivalencias")

# ---------------------------------------------------------------------
# 1. New "Trabajo de Grado" rows on PlanEstudios
# ---------------------------------------------------------------------

$ws1.Cells.Item(66, 1).Value = 2027633
$ws1.Cells.Item(66, 5).Value = "P"
$ws1.Cells.Item(66, 2).Value = "Trabajo de Grado - Trabajos Investigativos"
$ws1.Cells.Item(66, 3).Value = 8
$ws1.Cells.Item(66, 4).Value = "TRABAJO DE GRADO"

$ws1.Cells.Item(67, 1).Value = 2027634
$ws1.Cells.Item(67, 5).Value = "P"
$ws1.Cells.Item(67, 2).Value = "Trabajo de Grado - Asignaturas de Posgrado"
$ws1.Cells.Item(67, 3).Value = 8
$ws1.Cells.Item(67, 4).Value = "TRABAJO DE GRADO"

$ws1.Cells.Item(68, 1).Value = 2027636
$ws1.Cells.Item(68, 5).Value = "P"
$ws1.Cells.Item(68, 2).Value = "Trabajo de Grado - Pasantías"
$ws1.Cells.Item(68, 3).Value = 8
$ws1.Cells.Item(68, 4).Value = "TRABAJO DE GRADO"

# ---------------------------------------------------------------------
# 2. Agrupaciones: "Libre Elección" credit count 39 -> 28
# ---------------------------------------------------------------------

$ws2.Cells.Item(11, 4).Value = 28

# ---------------------------------------------------------------------
# 3. Strip the leftover explicit cell styles (all of them resolved to the
#    same, default "Normal" format, so Excel collapses the style table
#    down to a single cellXfs entry once nothing references the others)
# ---------------------------------------------------------------------

function Clear-CellStyles($ws, [string[]]$addrs) {
    foreach ($a in $addrs) {
        $ws.Range($a).Style = "Normal"
    }
}

Clear-CellStyles $ws1 @(
    "C2", "C3", "C4", "E4", "C5", "E5", "C6", "C7", "D7", "D8",
    "D9", "D10", "D11", "D12", "D13", "D14", "D15",
    "D17", "D18", "D19", "D20", "D21",
    "D23", "D24", "D25",
    "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34",
    "D39", "D40", "D41", "D42", "D43",
    "D45", "D46", "D47", "D48", "D49", "D50", "D51",
    "D53", "E53", "D54", "E54", "D55", "E55", "D56", "E56",
    "D57", "E57", "D58", "E58", "D59", "E59", "D60", "E60",
    "D61", "E61", "D62", "E62", "D63", "E63", "D64", "E64",
    "D65", "E65"
)

Clear-CellStyles $ws2 @("A2", "B4", "B6", "B7", "A8", "B8")

Clear-CellStyles $ws3 @(
    "C2", "D2", "C3", "D3", "C4", "D4", "C5", "D5",
    "C6", "D6", "C7", "D7", "C8", "D8"
)

# ---------------------------------------------------------------------
# 4. View state: PlanEstudios scrolled near the new rows (no longer the
#    tab in focus), Agrupaciones becomes the active/selected sheet
# ---------------------------------------------------------------------

$ws1.Range("B68").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("E11").Select() | Out-Null
